$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to text
# so Excel keeps them as strings (matching original inlineStr cells),
# not auto-converted to numeric cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.137'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.414'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000208'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.613'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000118'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.59'
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.92'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '171.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0856'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.897'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.47'
$ws.Range("D51").Style = "Normal"

# Remaining cells (already safe as text, e.g. multi-dot prices, padded
# percentages, coin names and links) can be assigned directly.
$ws.Range("D2").Value = '66.207.62'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '3.563.90'
$ws.Range("E3").Value = '  +4.90%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  +1.92%  '
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("D7").Value = '3.562.27'
$ws.Range("E7").Value = '  +4.89%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +5.83%  '
$ws.Range("E10").Value = '  +2.49%  '
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("E12").Value = '  +2.39%  '
$ws.Range("D13").Value = '4.165.36'
$ws.Range("E13").Value = '  +4.99%  '
$ws.Range("E14").Value = '  +4.56%  '
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '3.557.38'
$ws.Range("E16").Value = '  +5.08%  '
$ws.Range("D17").Value = '66.282.94'
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("E19").Value = '  +10.62%  '
$ws.Range("E20").Value = '  +1.86%  '
$ws.Range("E21").Value = '  +2.62%  '
$ws.Range("E22").Value = '  +4.41%  '
$ws.Range("E23").Value = '  +5.95%  '
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("D25").Value = '3.701.70'
$ws.Range("E25").Value = '  +4.80%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +9.43%  '
$ws.Range("E28").Value = '  +4.72%  '
$ws.Range("E29").Value = '  +3.42%  '
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("E33").Value = '  +5.03%  '
$ws.Range("D34").Value = '3.555.24'
$ws.Range("E34").Value = '  +4.85%  '
$ws.Range("E35").Value = '  -3.23%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  +4.85%  '
$ws.Range("E38").Value = '  +5.89%  '
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("E41").Value = '  +1.93%  '
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("E43").Value = '  +4.03%  '
$ws.Range("E44").Value = '  +3.32%  '
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("E46").Value = '  +1.80%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("E47").Value = '  +2.80%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("E49").Value = '  +5.01%  '
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("E51").Value = '  +17.09%  '
